# Applies per-cell text updates to match the Dec 23 2023 cryptos data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.744.17"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.312.21"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'108.43"
$ws.Range("E5").Value = "  +10.07%  "
$ws.Range("D6").Value = "'271.17"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'48.15"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("D11").Value = "'0.0940"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "'8.36"
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "2.654.39"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "'0.865"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "2.310.53"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "43.787.43"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "'6.33"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'72.36"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'2.52"
$ws.Range("E22").Value = "  +8.89%  "
$ws.Range("D23").Value = "'234.63"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("D24").Value = "'2.90"
$ws.Range("E24").Value = "  +14.35%  "
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'11.40"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "'41.94"
$ws.Range("E28").Value = "  +9.07%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "'177.69"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'21.98"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "'0.0919"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "'5.62"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +7.29%  "
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  +15.41%  "
$ws.Range("D40").Value = "'2.37"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").Value = "'67.47"
$ws.Range("E43").Value = "  +7.51%  "
$ws.Range("D44").Value = "'12.38"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "'5.59"
$ws.Range("E45").Value = "  +4.69%  "
$ws.Range("D46").Value = "'8.83"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'99.87"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.22"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "'0.438"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").Value = "2.546.10"
$ws.Range("E51").Value = "  -0.29%  "
